$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "inline" rule's description (row 2, column C) to include an
# escaped angle-bracket sequence.
$ws.Range("C2").Value = "A few words`ninside &gt; &amp;gt; inline"

# Insert a new row for the "infinity" norm rule right after the "inline" rule
# (row 2), pushing every following row down by one.
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "my-chapter_name"
$ws.Range("B3").Value = "infinity"
$ws.Range("C3").Value = "abc &amp;#x221e; def"
$ws.Range("D3").Value = "[`"norm:infinity`"]"

# Keep the worksheet Table in sync with the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F33"))
